$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells: "_old" -> "_FV2210", "_new" -> "_FV2304" ---
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = $cell.Value2
    if ($text -like "*_old") {
        $cell.Value = ($text -replace "_old$", "_FV2210")
    } elseif ($text -like "*_new") {
        $cell.Value = ($text -replace "_new$", "_FV2304")
    }
}

# --- 2. Freeze the header row (pane split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the data range into a table (adds autoFilter + tableParts) ---
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U90"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# --- 4. Restore the active selection to A1 ---
$ws.Range("A1").Select()
